# Auto-generated Excel COM-interop script
# Applies the Tonberry_Profits leve-profit recalculation update across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (ALC)
$ws.Range("H18").Value = 16248.75
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# Row 40 (ALC)
$ws.Range("H40").Value = 2265.9
$ws.Range("I40").Value = 2200
$ws.Range("J40").Value = 2331.8
$ws.Range("K40").Value = 2200
$ws.Range("L40").Value = 2331.8
$ws.Range("M40").Value = -2025
$ws.Range("N40").Value = -2681.8

# Row 47 (ALC)
$ws.Range("H47").Value = 9689.333000000001
$ws.Range("I47").Value = 7034
$ws.Range("J47").Value = 15000
$ws.Range("K47").Value = 7034
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = -6062
$ws.Range("N47").Value = -16944

# Row 61 (ALC)
$ws.Range("H61").Value = 3250
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 12000
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -11828
$ws.Range("N61").Value = -7844

# Row 86 (ALC)
$ws.Range("H86").Value = 2687
$ws.Range("J86").Value = 5500
$ws.Range("L86").Value = 5500
$ws.Range("N86").Value = -7746

# Row 89 (ALC)
$ws.Range("H89").Value = 2687
$ws.Range("J89").Value = 5500
$ws.Range("L89").Value = 27500
$ws.Range("N89").Value = -38732

# Row 104 (ALC)
$ws.Range("H104").Value = 2292
$ws.Range("I104").Value = 2820
$ws.Range("K104").Value = 8460
$ws.Range("M104").Value = -6713

# Row 112 (ALC)
$ws.Range("H112").Value = 4408.8335
$ws.Range("J112").Value = 4550.5293
$ws.Range("L112").Value = 13651.5879
$ws.Range("N112").Value = -15867.5879

# Row 132 (ALC)
$ws.Range("H132").Value = 1082.8718
$ws.Range("I132").Value = 984.6486
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 2953.9458
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -423.9458
$ws.Range("N132").Value = -13760

# Row 137 (ALC)
$ws.Range("H137").Value = 1884.238
$ws.Range("I137").Value = 1649.2
$ws.Range("J137").Value = 2097.9092
$ws.Range("K137").Value = 4947.6
$ws.Range("L137").Value = 6293.7276
$ws.Range("M137").Value = -2397.6
$ws.Range("N137").Value = -11393.7276

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1088.5
$ws.Range("I2").Value = 1155.5
$ws.Range("J2").Value = 887.5
$ws.Range("K2").Value = 1155.5
$ws.Range("L2").Value = 887.5
$ws.Range("M2").Value = -1042.5
$ws.Range("N2").Value = -1113.5

# Row 32 (ARM)
$ws.Range("H32").Value = 2576.9023
$ws.Range("I32").Value = 1769.9265
$ws.Range("K32").Value = 1769.9265
$ws.Range("M32").Value = -1482.9265

# Row 61 (ARM)
$ws.Range("H61").Value = 2178.4644
$ws.Range("I61").Value = 1118.3
$ws.Range("K61").Value = 1118.3
$ws.Range("M61").Value = -906.3

# Row 74 (ARM)
$ws.Range("H74").Value = 1658.8
$ws.Range("I74").Value = 1598.1111
$ws.Range("J74").Value = 1708.4546
$ws.Range("K74").Value = 1598.1111
$ws.Range("L74").Value = 1708.4546
$ws.Range("M74").Value = -724.1111000000001
$ws.Range("N74").Value = -3456.4546

# Row 77 (ARM)
$ws.Range("H77").Value = 1658.8
$ws.Range("I77").Value = 1598.1111
$ws.Range("J77").Value = 1708.4546
$ws.Range("K77").Value = 7990.5555
$ws.Range("L77").Value = 8542.273000000001
$ws.Range("M77").Value = -3622.5555
$ws.Range("N77").Value = -17278.273

# Row 110 (ARM)
$ws.Range("H110").Value = 1658
$ws.Range("I110").Value = 1015.3077
$ws.Range("K110").Value = 1015.3077
$ws.Range("M110").Value = 1029.6923

# Row 116 (ARM)
$ws.Range("H116").Value = 1088.5
$ws.Range("I116").Value = 1155.5
$ws.Range("J116").Value = 887.5
$ws.Range("K116").Value = 1155.5
$ws.Range("L116").Value = 887.5
$ws.Range("M116").Value = 1138.5
$ws.Range("N116").Value = -5475.5

# Row 122 (ARM)
$ws.Range("H122").Value = 599
$ws.Range("I122").Value = 599
$ws.Range("K122").Value = 1797
$ws.Range("M122").Value = 653

# Row 124 (ARM)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Row 125 (ARM)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132 (ARM)
$ws.Range("H132").Value = 1347.3392
$ws.Range("I132").Value = 1000.2955
$ws.Range("K132").Value = 3000.8865
$ws.Range("M132").Value = -470.8864999999996

# Row 136 (ARM)
$ws.Range("H136").Value = 2178.4644
$ws.Range("I136").Value = 1118.3
$ws.Range("K136").Value = 3354.9
$ws.Range("M136").Value = -804.8999999999996

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1088.5
$ws.Range("I3").Value = 1155.5
$ws.Range("J3").Value = 887.5
$ws.Range("K3").Value = 1155.5
$ws.Range("L3").Value = 887.5
$ws.Range("M3").Value = -1041.5
$ws.Range("N3").Value = -1115.5

# Row 105 (BSM)
$ws.Range("H105").Value = 2406.6191
$ws.Range("I105").Value = 2372.2942
$ws.Range("J105").Value = 2552.5
$ws.Range("K105").Value = 2372.2942
$ws.Range("L105").Value = 2552.5
$ws.Range("M105").Value = -625.2941999999998
$ws.Range("N105").Value = -6046.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1667.2727
$ws.Range("I31").Value = 1245.6
$ws.Range("J31").Value = 2570.8572
$ws.Range("K31").Value = 1245.6
$ws.Range("L31").Value = 2570.8572
$ws.Range("M31").Value = -950.5999999999999
$ws.Range("N31").Value = -3160.8572

# Row 34 (CRP)
$ws.Range("H34").Value = 1667.2727
$ws.Range("I34").Value = 1245.6
$ws.Range("J34").Value = 2570.8572
$ws.Range("K34").Value = 1245.6
$ws.Range("L34").Value = 2570.8572
$ws.Range("M34").Value = -1043.6
$ws.Range("N34").Value = -2974.8572

# Row 107 (CRP)
$ws.Range("H107").Value = 510.33334
$ws.Range("I107").Value = 427.1111
$ws.Range("J107").Value = 760
$ws.Range("K107").Value = 427.1111
$ws.Range("L107").Value = 760
$ws.Range("M107").Value = 1492.8889
$ws.Range("N107").Value = -4600

# Row 124 (CRP)
$ws.Range("H124").Value = 19400
$ws.Range("J124").Value = 19400
$ws.Range("L124").Value = 19400
$ws.Range("N124").Value = -24310

# Row 134 (CRP)
$ws.Range("H134").Value = 2222.318
$ws.Range("I134").Value = 1904.8422
$ws.Range("K134").Value = 5714.5266
$ws.Range("M134").Value = -3179.5266

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (CUL)
$ws.Range("H33").Value = 92.59999999999999
$ws.Range("J33").Value = 102
$ws.Range("L33").Value = 612
$ws.Range("N33").Value = -1178

# Row 74 (CUL)
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 6000
$ws.Range("M74").Value = -4939

# Row 77 (CUL)
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 18000
$ws.Range("M77").Value = -12696

# Row 131 (CUL)
$ws.Range("H131").Value = 2800.3
$ws.Range("J131").Value = 2940.1382
$ws.Range("L131").Value = 8820.4146
$ws.Range("N131").Value = -18900.4146

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 1402.375
$ws.Range("I97").Value = 1304.2727
$ws.Range("J97").Value = 1618.2
$ws.Range("K97").Value = 1304.2727
$ws.Range("L97").Value = 1618.2
$ws.Range("M97").Value = -808.2727
$ws.Range("N97").Value = -2610.2

# Row 132 (GSM)
$ws.Range("H132").Value = 2783.077
$ws.Range("I132").Value = 2537.8696
$ws.Range("J132").Value = 4663
$ws.Range("K132").Value = 7613.6088
$ws.Range("L132").Value = 13989
$ws.Range("M132").Value = -5083.6088
$ws.Range("N132").Value = -19049

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2901.5625
$ws.Range("I7").Value = 1954.6666
$ws.Range("K7").Value = 1954.6666
$ws.Range("M7").Value = -1842.6666

# Row 46 (LTW)
$ws.Range("H46").Value = 2992
$ws.Range("J46").Value = 3257.3333
$ws.Range("L46").Value = 3257.3333
$ws.Range("N46").Value = -3633.3333

# Row 122 (LTW)
$ws.Range("H122").Value = 5503.846
$ws.Range("I122").Value = 4850.6665
$ws.Range("K122").Value = 14551.9995
$ws.Range("M122").Value = -12101.9995

# Row 126 (LTW)
$ws.Range("H126").Value = 2901.5625
$ws.Range("I126").Value = 1954.6666
$ws.Range("K126").Value = 5863.9998
$ws.Range("M126").Value = -3393.9998

# Row 132 (LTW)
$ws.Range("H132").Value = 2552.0312
$ws.Range("I132").Value = 2291
$ws.Range("J132").Value = 2708.65
$ws.Range("K132").Value = 6873
$ws.Range("L132").Value = 8125.950000000001
$ws.Range("M132").Value = -4343
$ws.Range("N132").Value = -13185.95

# Row 136 (LTW)
$ws.Range("H136").Value = 3300.3809
$ws.Range("I136").Value = 2747.1
$ws.Range("J136").Value = 3803.3635
$ws.Range("K136").Value = 8241.299999999999
$ws.Range("L136").Value = 11410.0905
$ws.Range("M136").Value = -5691.299999999999
$ws.Range("N136").Value = -16510.0905
